$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are treated as text so numeric-looking strings
# (prices, percentages, hour codes) are preserved exactly as text, not
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "295.15"
$ws.Range("D3").Value = "31.16"
$ws.Range("D4").Value = "4.920"
$ws.Range("D5").Value = "0.07444"
$ws.Range("D6").Value = "2.281"
$ws.Range("D7").Value = "7.759"
$ws.Range("D9").Value = "0.9146"
$ws.Range("D10").Value = "0.09065"
$ws.Range("D11").Value = "0.1710"
$ws.Range("D12").Value = "0.08325"
$ws.Range("D13").Value = "0.03119"
$ws.Range("D14").Value = "0.1007"
$ws.Range("D15").Value = "0.001514"
$ws.Range("D16").Value = "0.005728"
$ws.Range("D17").Value = "3.503"
$ws.Range("D18").Value = "2.080"
$ws.Range("D19").Value = "0.3327"
$ws.Range("D21").Value = "3.980"
$ws.Range("D23").Value = "0.04548"
$ws.Range("D25").Value = "0.004620"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("D27").Value = "0.0003397"
$ws.Range("D39").Value = "0.01609"
$ws.Range("D40").Value = "0.04486"
$ws.Range("D41").Value = "0.007332"
$ws.Range("D42").Value = "0.008979"
$ws.Range("D44").Value = "0.001963"
$ws.Range("D45").Value = "0.008600"
$ws.Range("D46").Value = "0.00006032"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("D48").Value = "2.299"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("D51").Value = "0.0002002"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "1.51%"
$ws.Range("E3").Value = "1.21%"
$ws.Range("E4").Value = "-0.74%"
$ws.Range("E5").Value = "3.98%"
$ws.Range("E6").Value = "27.99%"
$ws.Range("E7").Value = "1.29%"
$ws.Range("E8").Value = "0.29%"
$ws.Range("E9").Value = "2.18%"
$ws.Range("E10").Value = "17.85%"
$ws.Range("E11").Value = "3.84%"
$ws.Range("E12").Value = "3.97%"
$ws.Range("E13").Value = "3.05%"
$ws.Range("E14").Value = "0.51%"
$ws.Range("E15").Value = "0.72%"
$ws.Range("E16").Value = "-0.09%"
$ws.Range("E17").Value = "0.96%"
$ws.Range("E18").Value = "-0.13%"
$ws.Range("E19").Value = "1.48%"
$ws.Range("E20").Value = "-0.04%"
$ws.Range("E21").Value = "-1.59%"
$ws.Range("E23").Value = "0.73%"
$ws.Range("E24").Value = "0.12%"
$ws.Range("E25").Value = "15.39%"
$ws.Range("E26").Value = "4.07%"
$ws.Range("E39").Value = "0.31%"
$ws.Range("E40").Value = "2.78%"
$ws.Range("E41").Value = "-0.14%"
$ws.Range("E43").Value = "1.77%"
$ws.Range("E44").Value = "-4.16%"
$ws.Range("E45").Value = "-6.84%"
$ws.Range("E46").Value = "5.66%"
$ws.Range("E47").Value = "0.12%"
$ws.Range("E48").Value = "2.30%"
$ws.Range("E49").Value = "-33.25%"
$ws.Range("E50").Value = "0.12%"
$ws.Range("E51").Value = "0.12%"

# --- Column G (Hora) updates: all rows 2-51 go from 3 to 4 ---
$ws.Range("G2:G51").Value = "4"

# Restore default (Normal) style on touched ranges so no stray
# number-format style is left applied to the cells.
$ws.Range("D2:D51").Style = "Normal"
$ws.Range("E2:E51").Style = "Normal"
$ws.Range("G2:G51").Style = "Normal"
